$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.188994288444519
$ws.Range("B1").Value = 2.164887189865112
$ws.Range("C1").Value = 4.521456718444824
$ws.Range("D1").Value = 2.77139139175415
$ws.Range("E1").Value = 1.214289546012878
